# Generate Report for Handback
#
# Adds a new tracked file "b24b80c8-9f5e-443e-a787-66871520e937.md" as row 39
# to all three worksheets (Overview, zh-cn, de-de), mirroring the existing
# pattern used for the previous last row (row 38, file
# "a104e4c3-f3ac-4fd5-8af3-8c55d95d32bd").

$wb = $excel.ActiveWorkbook

$newFile = "b24b80c8-9f5e-443e-a787-66871520e937.md"
$newHash = "373fd05a24aea35216201b8465d9899636e2fc1f"
$zhXlf   = "b24b80c8-9f5e-443e-a787-66871520e937.373fd05a24aea35216201b8465d9899636e2fc1f.zh-cn.xlf"
$deXlf   = "b24b80c8-9f5e-443e-a787-66871520e937.373fd05a24aea35216201b8465d9899636e2fc1f.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"
$include = "Include"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> row 39 : A (File Name, hyperlink), B (zh-cn status),
# C (de-de status)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$aCell = $wsOverview.Range("A39")
$wsOverview.Hyperlinks.Add($aCell, "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$newFile", "", "", $newFile) | Out-Null

$wsOverview.Range("B39").Value = $statusInSync
$wsOverview.Range("C39").Value = $statusInSync

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> row 39
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$aCellZh = $wsZh.Range("A39")
$wsZh.Hyperlinks.Add($aCellZh, "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$newFile", "", "", $newFile) | Out-Null

$wsZh.Range("B39").Value = $statusInSync

$cCellZh = $wsZh.Range("C39")
$wsZh.Hyperlinks.Add($cCellZh, "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf", "", "", $zhXlf) | Out-Null

$wsZh.Range("D39").Value = "2016-03-03 09:42:27"

$eCellZh = $wsZh.Range("E39")
$wsZh.Hyperlinks.Add($eCellZh, "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0000000000000000000000000000000000000000/e2e/$newFile", "", "", $newFile) | Out-Null

$fCellZh = $wsZh.Range("F39")
$wsZh.Hyperlinks.Add($fCellZh, "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf", "", "", $zhXlf) | Out-Null

$wsZh.Range("G39").Value = "2016-03-03 09:43:14"
$wsZh.Range("H39").Value = $include
$wsZh.Range("I39").Value = ""

# ---------------------------------------------------------------------------
# Sheet "de-de" -> row 39
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$aCellDe = $wsDe.Range("A39")
$wsDe.Hyperlinks.Add($aCellDe, "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$newFile", "", "", $newFile) | Out-Null

$wsDe.Range("B39").Value = $statusInSync

$cCellDe = $wsDe.Range("C39")
$wsDe.Hyperlinks.Add($cCellDe, "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf", "", "", $deXlf) | Out-Null

$wsDe.Range("D39").Value = "2016-03-03 09:42:40"

$eCellDe = $wsDe.Range("E39")
$wsDe.Hyperlinks.Add($eCellDe, "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0000000000000000000000000000000000000000/e2e/$newFile", "", "", $newFile) | Out-Null

$fCellDe = $wsDe.Range("F39")
$wsDe.Hyperlinks.Add($fCellDe, "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf", "", "", $deXlf) | Out-Null

$wsDe.Range("G39").Value = "2016-03-03 09:43:35"
$wsDe.Range("H39").Value = $include
$wsDe.Range("I39").Value = ""

Write-Host "Row 39 added to Overview, zh-cn, de-de"
